$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: Dietabs price changes from 5.99 to 5
$ws.Range("C51").Value = 5

# Row 52 (Allergy Eye Drops) row height goes from 18.75 -> 18 (match the other "18" rows)
$ws.Range("A52:C52").RowHeight = 18

# New row 53: Neozep
$ws.Range("A52:C52").Copy()
$ws.Range("A53:C53").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = "Neozep"
$ws.Range("C53").Value = 1.09
$ws.Range("A53:C53").RowHeight = 18

# New row 54: blank placeholder row (quote-prefixed blank cell in B54)
$ws.Range("A51:C51").Copy()
$ws.Range("A54:C54").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A54").Value = ""
$ws.Range("B54").Value = "'"
$ws.Range("C54").Value = ""
$ws.Range("A54:C54").RowHeight = 18.75
